# Exam Center Final commit
# Updates the "Stage" and "Prod" sheets with the final exam-center time
# slot values and syncs the active selection on both sheets to A2:XFD2.

$wb = $excel.ActiveWorkbook

# --- Stage sheet ---
$wsStage = $wb.Worksheets.Item("Stage")
$wsStage.Range("C2").Value = 31
$wsStage.Range("D2").Value = 17
$wsStage.Range("A2:XFD2").Select() | Out-Null

# --- Prod sheet ---
$wsProd = $wb.Worksheets.Item("Prod")
$wsProd.Range("C2").Value = 31
$wsProd.Range("D2").Value = 17
$wsProd.Range("E2").Value = 59
$wsProd.Range("F2").Value = 22
$wsProd.Range("G2").Value = 57
$wsProd.Range("A2:XFD2").Select() | Out-Null
